$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'25.900.58"
$ws.Range("E2").Value2 = "  +0.22%  "
$ws.Range("D3").Value2 = "'1.740.55"
$ws.Range("E3").Value2 = "  +0.30%  "
$ws.Range("D4").Value2 = "'0.9983"
$ws.Range("E4").Value2 = "  -0.13%  "
$ws.Range("D5").Value2 = "'241.60"
$ws.Range("E5").Value2 = "  +5.31%  "
$ws.Range("D6").Value2 = "'0.9998"
$ws.Range("E6").Value2 = "  +0.04%  "
$ws.Range("D7").Value2 = "'0.5196"
$ws.Range("E7").Value2 = "  -1.28%  "
$ws.Range("D8").Value2 = "'0.2757"
$ws.Range("E8").Value2 = "  +0.29%  "
$ws.Range("D9").Value2 = "'0.06165"
$ws.Range("E9").Value2 = "  +0.39%  "
$ws.Range("D10").Value2 = "'1.742.54"
$ws.Range("E10").Value2 = "  +0.27%  "
$ws.Range("D11").Value2 = "'0.07188"
$ws.Range("E11").Value2 = "  +1.62%  "
$ws.Range("D12").Value2 = "'15.01"
$ws.Range("E12").Value2 = "  +0.30%  "
$ws.Range("D13").Value2 = "'0.6450"
$ws.Range("E13").Value2 = "  +0.75%  "
$ws.Range("D14").Value2 = "'4.611"
$ws.Range("E14").Value2 = "  +2.00%  "
$ws.Range("D15").Value2 = "'77.66"
$ws.Range("E15").Value2 = "  +1.17%  "
$ws.Range("D16").Value2 = "'1.000"
$ws.Range("E16").Value2 = "  +0.12%  "
$ws.Range("D17").Value2 = "'0.9983"
$ws.Range("E17").Value2 = "  -0.14%  "
$ws.Range("D18").Value2 = "'25.908.14"
$ws.Range("E18").Value2 = "  +0.32%  "
$ws.Range("D19").Value2 = "'11.73"
$ws.Range("E19").Value2 = "  +1.91%  "
$ws.Range("D20").Value2 = "'0.000006770"
$ws.Range("E20").Value2 = "  +1.77%  "
$ws.Range("D21").Value2 = "'1.965.15"
$ws.Range("E21").Value2 = "  +0.31%  "
$ws.Range("D22").Value2 = "'4.285"
$ws.Range("E22").Value2 = "  +1.02%  "
$ws.Range("D23").Value2 = "'8.636"
$ws.Range("E23").Value2 = "  -1.56%  "
$ws.Range("D24").Value2 = "'5.271"
$ws.Range("E24").Value2 = "  +2.21%  "
$ws.Range("D25").Value2 = "'139.02"
$ws.Range("E25").Value2 = "  -0.77%  "
$ws.Range("D26").Value2 = "'1.517"
$ws.Range("E26").Value2 = "  +0.35%  "
$ws.Range("D27").Value2 = "'15.19"
$ws.Range("E27").Value2 = "  +0.58%  "
$ws.Range("D28").Value2 = "'1.767"
$ws.Range("E28").Value2 = "  -0.75%  "
$ws.Range("D29").Value2 = "'106.12"
$ws.Range("E29").Value2 = "  +3.94%  "
$ws.Range("D30").Value2 = "'3.925"
$ws.Range("E30").Value2 = "  +5.74%  "
$ws.Range("D31").Value2 = "'0.08289"
$ws.Range("E31").Value2 = "  -0.41%  "
$ws.Range("D32").Value2 = "'3.707"
$ws.Range("E32").Value2 = "  +5.16%  "
$ws.Range("D33").Value2 = "'0.04618"
$ws.Range("E33").Value2 = "  +3.07%  "
$ws.Range("D34").Value2 = "'2.643"
$ws.Range("E34").Value2 = "  +1.22%  "
$ws.Range("D35").Value2 = "'0.9899"
$ws.Range("E35").Value2 = "  +1.90%  "
$ws.Range("D36").Value2 = "'0.6188"
$ws.Range("E36").Value2 = "  +0.05%  "
$ws.Range("D37").Value2 = "'2.677"
$ws.Range("E37").Value2 = "  +0.10%  "
$ws.Range("D38").Value2 = "'0.01609"
$ws.Range("E38").Value2 = "  +2.48%  "
$ws.Range("D39").Value2 = "'1.933"
$ws.Range("E39").Value2 = "  +1.81%  "
$ws.Range("D40").Value2 = "'0.9997"
$ws.Range("E40").Value2 = "  +0.10%  "
$ws.Range("D41").Value2 = "'97.64"
$ws.Range("E41").Value2 = "  -2.35%  "
$ws.Range("D42").Value2 = "'0.3855"
$ws.Range("E42").Value2 = "  +0.26%  "
$ws.Range("D43").Value2 = "'0.7419"
$ws.Range("E43").Value2 = "  +2.40%  "
$ws.Range("D44").Value2 = "'4.983"
$ws.Range("E44").Value2 = "  -0.84%  "
$ws.Range("D45").Value2 = "'0.1132"
$ws.Range("E45").Value2 = "  +1.25%  "
$ws.Range("D46").Value2 = "'6.250"
$ws.Range("E46").Value2 = "  +0.97%  "
$ws.Range("D47").Value2 = "'0.05239"
$ws.Range("E47").Value2 = "  -1.56%  "
$ws.Range("D48").Value2 = "'54.83"
$ws.Range("E48").Value2 = "  +2.84%  "
$ws.Range("D49").Value2 = "'30.45"
$ws.Range("E49").Value2 = "  +1.58%  "
$ws.Range("D50").Value2 = "'7.596"
$ws.Range("E50").Value2 = "  -0.45%  "
$ws.Range("D51").Value2 = "'0.3420"
$ws.Range("E51").Value2 = "  +0.78%  "

$ws.Range("D2:D51").Style = "Normal"

